$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Split the J1:J4 merged note block into J1:J2 (new github-help note)
#     and J3:J4 (old "Program does not read..." note) ---
$ws.Range("J1:J4").UnMerge()
$ws.Range("J1").Value = "Learn to use at https://github.com/MKSanic/myob-entry/tree/main"
$ws.Range("J3").Value = "Program does not read anything beyond column I"
$ws.Range("J1:J2").Merge()
$ws.Range("J3:J4").Merge()

# --- "Commonly used codes" box gains a second (label) column ---
# Give the numeric codes in M2:M7 their own label column (N2:N7),
# and replace M2:M7 with the new text labels.
$ws.Range("N2").Value = $ws.Range("M2").Value
$ws.Range("M2").Value = "Renewal"

$ws.Range("N3").Value = $ws.Range("M3").Value
$ws.Range("M3").Value = "Fund Base"

$ws.Range("N4").Value = $ws.Range("M4").Value
$ws.Range("M4").Value = "FAFs"

$ws.Range("N5").Value = $ws.Range("M5").Value
$ws.Range("M5").Value = "OR FAFs"

$ws.Range("N6").Value = $ws.Range("M6").Value
$ws.Range("M6").Value = "Initial Comm"

$ws.Range("N7").Value = $ws.Range("M7").Value
$ws.Range("M7").Value = "Override comm"

# --- Header cell M1 ("Commonly used codes") now spans M1:N1, centered, bold ---
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("M1:N1").Merge()

# --- Column J a touch wider to fit the new note text ---
$ws.Columns("J").ColumnWidth = 28.86

# --- Selection / scroll position ---
$ws.Range("K24").Select()
